$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 64.375
$ws.Range("I33").Value = 65
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 65
$ws.Range("L33").Value = 60
$ws.Range("M33").Value = 164
$ws.Range("N33").Value = -518
$ws.Range("H40").Value = 1760.3077
$ws.Range("I40").Value = 1731.7778
$ws.Range("J40").Value = 1824.5
$ws.Range("K40").Value = 1731.7778
$ws.Range("L40").Value = 1824.5
$ws.Range("M40").Value = -1556.7778
$ws.Range("N40").Value = -2174.5
$ws.Range("H53").Value = 337.33334
$ws.Range("I53").Value = 289.66666
$ws.Range("K53").Value = 289.66666
$ws.Range("M53").Value = 347.33334
$ws.Range("H86").Value = 6555.6
$ws.Range("I86").Value = 7445
$ws.Range("J86").Value = 5962.6665
$ws.Range("K86").Value = 7445
$ws.Range("L86").Value = 5962.6665
$ws.Range("M86").Value = -6322
$ws.Range("N86").Value = -8208.666499999999
$ws.Range("H89").Value = 6555.6
$ws.Range("I89").Value = 7445
$ws.Range("J89").Value = 5962.6665
$ws.Range("K89").Value = 37225
$ws.Range("L89").Value = 29813.3325
$ws.Range("M89").Value = -31609
$ws.Range("N89").Value = -41045.3325
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 5000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -6082
$ws.Range("H101").Value = 50000124
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 3929.6667
$ws.Range("I113").Value = 3644.5
$ws.Range("K113").Value = 3644.5
$ws.Range("M113").Value = -390.5
$ws.Range("H127").Value = 2933.5
$ws.Range("J127").Value = 1249
$ws.Range("L127").Value = 3747
$ws.Range("N127").Value = -13667
$ws.Range("H129").Value = 2668.8667
$ws.Range("I129").Value = 2317.25
$ws.Range("K129").Value = 6951.75
$ws.Range("M129").Value = -1951.75
$ws.Range("H131").Value = 21235.8
$ws.Range("I131").Value = 21235.8
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 63707.39999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -58667.39999999999
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2351.5854
$ws.Range("I132").Value = 1571.8572
$ws.Range("K132").Value = 4715.571599999999
$ws.Range("M132").Value = -2185.571599999999
$ws.Range("H135").Value = 87.666664
$ws.Range("I135").Value = 87.666664
$ws.Range("K135").Value = 788.9999759999999
$ws.Range("M135").Value = 1746.000024
$ws.Range("H141").Value = 4399
$ws.Range("I141").Value = 5332
$ws.Range("J141").Value = 2999.5
$ws.Range("K141").Value = 15996
$ws.Range("L141").Value = 8998.5
$ws.Range("M141").Value = -10816
$ws.Range("N141").Value = -19358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1291.1
$ws.Range("I61").Value = 1362.125
$ws.Range("J61").Value = 1007
$ws.Range("K61").Value = 1362.125
$ws.Range("L61").Value = 1007
$ws.Range("M61").Value = -1150.125
$ws.Range("N61").Value = -1431
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 1291.1
$ws.Range("I136").Value = 1362.125
$ws.Range("J136").Value = 1007
$ws.Range("K136").Value = 4086.375
$ws.Range("L136").Value = 3021
$ws.Range("M136").Value = -1536.375
$ws.Range("N136").Value = -8121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1654.3334
$ws.Range("I20").Value = 1654.3334
$ws.Range("K20").Value = 1654.3334
$ws.Range("M20").Value = -1407.3334
$ws.Range("H86").Value = 1320.375
$ws.Range("I86").Value = 1217.25
$ws.Range("K86").Value = 1217.25
$ws.Range("M86").Value = -94.25
$ws.Range("H89").Value = 1320.375
$ws.Range("I89").Value = 1217.25
$ws.Range("K89").Value = 6086.25
$ws.Range("M89").Value = -470.25
$ws.Range("H94").Value = 774.7368
$ws.Range("I94").Value = 769.6923
$ws.Range("J94").Value = 785.6667
$ws.Range("K94").Value = 769.6923
$ws.Range("L94").Value = 785.6667
$ws.Range("M94").Value = -318.6923
$ws.Range("N94").Value = -1687.6667
$ws.Range("H99").Value = 1195.5294
$ws.Range("I99").Value = 940.8
$ws.Range("K99").Value = 940.8
$ws.Range("M99").Value = 557.2
$ws.Range("H105").Value = 3061.625
$ws.Range("I105").Value = 3070.8572
$ws.Range("J105").Value = 2997
$ws.Range("K105").Value = 3070.8572
$ws.Range("L105").Value = 2997
$ws.Range("M105").Value = -1323.8572
$ws.Range("N105").Value = -6491
$ws.Range("H107").Value = 5007.25
$ws.Range("I107").Value = 3548.75
$ws.Range("K107").Value = 3548.75
$ws.Range("M107").Value = -1628.75
$ws.Range("H134").Value = 2529.625
$ws.Range("I134").Value = 2529.625
$ws.Range("K134").Value = 7588.875
$ws.Range("M134").Value = -5053.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2148.2856
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2148.2856
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 2148.2856
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -2738.2856
$ws.Range("H34").Value = 2148.2856
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2148.2856
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2148.2856
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -2552.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1242.5
$ws.Range("J5").Value = 1642
$ws.Range("L5").Value = 4926
$ws.Range("N5").Value = -5150
$ws.Range("H11").Value = 28866642
$ws.Range("I11").Value = 41696090
$ws.Range("J11").Value = 386.25
$ws.Range("K11").Value = 125088270
$ws.Range("L11").Value = 1158.75
$ws.Range("M11").Value = -125088130
$ws.Range("N11").Value = -1438.75
$ws.Range("H26").Value = 83
$ws.Range("I26").Value = 79
$ws.Range("J26").Value = 95
$ws.Range("K26").Value = 237
$ws.Range("L26").Value = 285
$ws.Range("M26").Value = 51
$ws.Range("N26").Value = -861
$ws.Range("H75").Value = 10071.167
$ws.Range("J75").Value = 18338.334
$ws.Range("L75").Value = 55015.00199999999
$ws.Range("N75").Value = -57011.00199999999
$ws.Range("H78").Value = 10071.167
$ws.Range("J78").Value = 18338.334
$ws.Range("L78").Value = 165045.006
$ws.Range("N78").Value = -175029.006
$ws.Range("H108").Value = 2679.5
$ws.Range("I108").Value = 2679.5
$ws.Range("K108").Value = 8038.5
$ws.Range("M108").Value = -5158.5
$ws.Range("H129").Value = 670094.4
$ws.Range("I129").Value = 2326.5715
$ws.Range("J129").Value = 1254391.2
$ws.Range("K129").Value = 6979.7145
$ws.Range("L129").Value = 3763173.6
$ws.Range("M129").Value = -1979.7145
$ws.Range("N129").Value = -3773173.6
$ws.Range("H131").Value = 528236.3
$ws.Range("I131").Value = 976.5
$ws.Range("J131").Value = 771587
$ws.Range("K131").Value = 2929.5
$ws.Range("L131").Value = 2314761
$ws.Range("M131").Value = 2110.5
$ws.Range("N131").Value = -2324841
$ws.Range("H135").Value = 1242.5
$ws.Range("J135").Value = 1642
$ws.Range("L135").Value = 14778
$ws.Range("N135").Value = -19848
$ws.Range("H139").Value = 998.5
$ws.Range("I139").Value = 998.5
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2995.5
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2144.5
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 10369.941
$ws.Range("I140").Value = 1213.4286
$ws.Range("J140").Value = 16779.5
$ws.Range("K140").Value = 3640.2858
$ws.Range("L140").Value = 50338.5
$ws.Range("M140").Value = 1539.7142
$ws.Range("N140").Value = -60698.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5500
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 5000
$ws.Range("M102").Value = -3378
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 58653.223
$ws.Range("J46").Value = 4500
$ws.Range("L46").Value = 4500
$ws.Range("N46").Value = -4876
$ws.Range("H61").Value = 3478.2
$ws.Range("I61").Value = 2848.375
$ws.Range("K61").Value = 2848.375
$ws.Range("M61").Value = -2646.375
$ws.Range("H113").Value = 3478.2
$ws.Range("I113").Value = 2848.375
$ws.Range("K113").Value = 2848.375
$ws.Range("M113").Value = -678.375
$ws.Range("H132").Value = 18628.334
$ws.Range("I132").Value = 26192.5
$ws.Range("K132").Value = 78577.5
$ws.Range("M132").Value = -76047.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8573.076999999999
$ws.Range("J136").Value = 4352.5
$ws.Range("L136").Value = 13057.5
$ws.Range("N136").Value = -18157.5
